$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Use API Output as Input for Logstash".
# A new bullet ("Or compare years?", nested one level deeper) was inserted directly
# above it, which pushes the existing "Use API Output..." / "Create configuration
# files..." / "Use Kibana..." / "Text analysis?" bullets down by one slot.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Use API Output as Input for Logstash") {
        $target = $p
        break
    }
}

# Insert a brand new (empty) paragraph right before it; Word clones the
# paragraph/run formatting (style, numbering, language) from $target. After
# the insertion, $target itself re-seats onto that new empty paragraph
# (the old text stays on the paragraph that follows it).
$target.Range.InsertParagraphBefore()
$newPara = $target

# Fill in the new bullet's text. A trailing placeholder character is typed
# first so the bookmark below can be anchored *after* the real text without
# landing on the paragraph-mark boundary (which re-seats to the paragraph
# boundary instead of trailing the run); the placeholder is then deleted.
$newPara.Range.Text = "Or compare years?X"

# Demote it one level (ilvl 0 -> 1, i.e. ListLevelNumber 1 -> 2).
$newPara.Range.ListFormat.ListLevelNumber = 2

# The "_GoBack" bookmark used to sit at the end of the old last bullet
# ("Text analysis?"). Move it: drop it from there ...
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ... and re-create it around the placeholder character, trailing the real
# text of the newly typed bullet, since that's where the author's cursor
# ended up after typing.
$full = $newPara.Range
$placeholder = $full.Duplicate
$placeholder.SetRange($full.End - 2, $full.End - 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

# Remove the placeholder character; the (now empty) bookmark stays put,
# trailing the run text.
$placeholder = $d.Range($full.End - 2, $full.End - 1)
$placeholder.Delete()
